$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D1 header cell: match the style already used by C1 (bold, centered).
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill D2:D8 with the "Last Edited" date (2020-09-17), matching the
# "Created Date" column's centered date style (copy format from C2, which
# already carries the date number format + centered alignment).
$ws.Range("C2").Copy()
$ws.Range("D2:D8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$lastEdited = Get-Date -Year 2020 -Month 9 -Day 17 -Hour 0 -Minute 0 -Second 0
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 4).Value = $lastEdited
}

# Clear any active-cell selection recorded in the sheet view.
$ws.Range("A1").Select()
